# Commit: rename GPIB address placeholder to a SOCKET address in the
# "Connections" lookup sheet of the equipment-connections workbook.
# (the GPIB::4 string used for a socket/ethernet example record is
# replaced by a fully-qualified SOCKET::host::port address)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connections")

# Make "Connections" the active sheet (it already is) and put the
# selection on the cell being edited, matching the saved UI state.
$ws.Activate()
$ws.Range("E3").Select()

# Update the cell's content: GPIB::4 -> SOCKET::1.2.3.4::1234
$ws.Range("E3").Value = "SOCKET::1.2.3.4::1234"
